$wb = $excel.ActiveWorkbook

# Template sheet for the new country sheets: "Germany" has the exact same
# column widths / row layout (19 rows) the new sheets are based on.
$template = $wb.Worksheets.Item("Germany")

# ---------------------------------------------------------------------
# Netherlands (19 rows, full product list incl. P32AR/P32DR)
# ---------------------------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsNL = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNL.Name = "Netherlands"

# set B4 before B2 so the shared-string table picks up the same order
# the original authoring tool used
$wsNL.Range("B4").Value2 = "NGC-3144/T2179"
$wsNL.Range("B4").Borders.LineStyle = -4142
$wsNL.Range("B2").Value2 = "Netherlands Market"

$wsNL.Activate()
$wsNL.Range("B4").Select()

# ---------------------------------------------------------------------
# Austria (17 rows -- P32AR/P32DR rows removed)
# ---------------------------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsAT = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAT.Name = "Austria"
$wsAT.Rows("16:17").Delete()

$wsAT.Range("B4").Value2 = "NGC-3817/T2275"
$wsAT.Range("B4").Borders.LineStyle = -4142
$wsAT.Range("B2").Value2 = "Austria Market"

$wsAT.Activate()
$wsAT.Range("B4").Select()

# ---------------------------------------------------------------------
# Denmark (18 rows -- P32DR row removed, P32AR row repurposed)
# ---------------------------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsDK = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsDK.Name = "Denmark"
$wsDK.Rows("17:17").Delete()

$wsDK.Range("A16").Value2 = "MZXSDR240"
$wsDK.Range("B4").Value2 = "NGC-2913/T2754"
$wsDK.Range("B4").Borders.LineStyle = -4142
$wsDK.Range("B2").Value2 = "Denmark Market"

$wsDK.Activate()
$wsDK.Range("B4").Select()

# ---------------------------------------------------------------------
# Netherlands ends up as the active tab (activeTab index 10)
# ---------------------------------------------------------------------
$wsNL.Activate()
$wsNL.Range("B4").Select()
